# Apply "correction de numero d'employe pour les administrateurs" edit
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# 1. Fill in "Responsable" column (F) for each task row
$ws.Range("F3").Value = "Serge"
$ws.Range("F4").Value = "Serge"
$ws.Range("F5").Value = "Pat"
$ws.Range("F6").Value = "Pat"
$ws.Range("F7").Value = "Pat"
$ws.Range("F8").Value = "Pat"
$ws.Range("F9").Value = "Serge"
$ws.Range("F10").Value = "Serge, Pat"
$ws.Range("F11").Value = "Pat"
$ws.Range("F13").Value = "Pat"
$ws.Range("F14").Value = "Pat"
$ws.Range("F15").Value = "Pat"
$ws.Range("F16").Value = "pat"
$ws.Range("F17").Value = "Pat"
$ws.Range("F18").Value = "Pat"
$ws.Range("F19").Value = "Serge"
$ws.Range("F21").Value = "Serge"

# 2. Fix task description text in B21 (ligne de commande -> GUI)
$ws.Range("B21").Value = "Ajouter une interface GUI pour faciliter l'utilisation"

# 3. Update the view: scroll position / selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H20").Select()
